# Update "想去人数" (F column) values for rows that changed between scrapes.
# Applies to both the "展览" and "全部类型" worksheets (they mirror the same data).

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3215
    4  = 233
    5  = 131
    6  = 202
    7  = 1694
    8  = 1635
    9  = 469
    10 = 373
    11 = 235
    21 = 56
    23 = 380
    24 = 215
    26 = 34
    29 = 287
    30 = 2194
    34 = 333
    35 = 570
    36 = 426
    40 = 519
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
